$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values can look numeric (e.g. "51.89", "37.403.06"); Excel
# would otherwise silently coerce them into numbers (losing the original text
# formatting, e.g. "1.00" -> 1). Force those cells to Text format before
# assignment, then clear the format again afterwards so the cell keeps the
# same (unstyled) appearance as in the original workbook.
$priceCells = @("D2", "D3", "D5", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D26", "D27", "D28", "D33", "D34", "D35", "D39", "D40", "D42", "D44", "D45", "D47", "D48", "D49", "D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '37.403.06'
$ws.Range("E2").Value = '  +2.99%  '
$ws.Range("D3").Value = '2.096.78'
$ws.Range("E3").Value = '  +4.42%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").Value = '250.71'
$ws.Range("E5").Value = '  +2.22%  '
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").Value = '51.89'
$ws.Range("E8").Value = '  +15.49%  '
$ws.Range("D9").Value = '61.38'
$ws.Range("E9").Value = '  +9.21%  '
$ws.Range("D10").Value = '0.375'
$ws.Range("E10").Value = '  +3.92%  '
$ws.Range("D11").Value = '0.0744'
$ws.Range("E11").Value = '  +3.71%  '
$ws.Range("E12").Value = '  +6.79%  '
$ws.Range("D13").Value = '15.24'
$ws.Range("E13").Value = '  +5.51%  '
$ws.Range("D14").Value = '2.401.92'
$ws.Range("E14").Value = '  +4.70%  '
$ws.Range("D15").Value = '0.830'
$ws.Range("E15").Value = '  +4.40%  '
$ws.Range("D16").Value = '2.103.63'
$ws.Range("E16").Value = '  +4.44%  '
$ws.Range("D17").Value = '5.10'
$ws.Range("E17").Value = '  +4.80%  '
$ws.Range("D18").Value = '37.334.71'
$ws.Range("E18").Value = '  +2.98%  '
$ws.Range("D19").Value = '72.02'
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '13.95'
$ws.Range("E20").Value = '  +8.37%  '
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0833'
$ws.Range("E21").Value = '  +2.73%  '
$ws.Range("D22").Value = '240.09'
$ws.Range("E22").Value = '  +2.91%  '
$ws.Range("E23").Value = '  +4.80%  '
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("D26").Value = '170.27'
$ws.Range("E26").Value = '  +5.47%  '
$ws.Range("D27").Value = '9.16'
$ws.Range("E27").Value = '  +8.79%  '
$ws.Range("D28").Value = '20.69'
$ws.Range("E28").Value = '  +4.82%  '
$ws.Range("E29").Value = '  +0.90%  '
$ws.Range("E30").Value = '  +0.55%  '
$ws.Range("E31").Value = '  +26.88%  '
$ws.Range("E32").Value = '  +2.99%  '
$ws.Range("D33").Value = '0.0609'
$ws.Range("E33").Value = '  +5.10%  '
$ws.Range("D34").Value = '0.0916'
$ws.Range("E34").Value = '  +8.86%  '
$ws.Range("D35").Value = '19.81'
$ws.Range("E35").Value = '  -7.33%  '
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("E37").Value = '  +9.07%  '
$ws.Range("E38").Value = '  -0.82%  '
$ws.Range("D39").Value = '4.09'
$ws.Range("E39").Value = '  +2.18%  '
$ws.Range("D40").Value = '1.32'
$ws.Range("E40").Value = '  -1.50%  '
$ws.Range("E41").Value = '  +11.00%  '
$ws.Range("D42").Value = '0.0224'
$ws.Range("E42").Value = '  +4.67%  '
$ws.Range("E43").Value = '  +8.83%  '
$ws.Range("D44").Value = '98.94'
$ws.Range("E44").Value = '  +2.24%  '
$ws.Range("D45").Value = '0.0897'
$ws.Range("E45").Value = '  +11.14%  '
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("D47").Value = '3.03'
$ws.Range("E47").Value = '  +9.46%  '
$ws.Range("D48").Value = '1.319.72'
$ws.Range("E48").Value = '  +1.17%  '
$ws.Range("D49").Value = '6.98'
$ws.Range("E49").Value = '  +15.53%  '
$ws.Range("D50").Value = '2.282.60'
$ws.Range("E50").Value = '  +4.83%  '
$ws.Range("E51").Value = '  +3.43%  '

foreach ($addr in $priceCells) {
    $ws.Range($addr).ClearFormats()
}
